$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.037.99'
$ws.Range("E2").Value = '  +2.35%  '

$ws.Range("D3").Value = '2.460.34'
$ws.Range("E3").Value = '  +2.01%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.28'
$ws.Range("E5").Value = '  +1.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.72'
$ws.Range("E6").Value = '  +2.40%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.542'
$ws.Range("E8").Value = '  +0.95%  '

$ws.Range("D9").Value = '2.458.70'
$ws.Range("E9").Value = '  +1.45%  '

$ws.Range("E10").Value = '  +2.64%  '

$ws.Range("E11").Value = '  +1.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.28'
$ws.Range("E12").Value = '  +1.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("E13").Value = '  +2.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.84'
$ws.Range("E14").Value = '  +9.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000179'
$ws.Range("E15").Value = '  +3.30%  '

$ws.Range("D16").Value = '2.906.90'
$ws.Range("E16").Value = '  +2.06%  '

$ws.Range("D17").Value = '62.928.29'

$ws.Range("D18").Value = '2.459.01'
$ws.Range("E18").Value = '  +1.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.03'
$ws.Range("E19").Value = '  +0.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.11'
$ws.Range("E20").Value = '  +4.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.35'
$ws.Range("E21").Value = '  +1.96%  '

$ws.Range("E22").Value = '  +12.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.14'
$ws.Range("E23").Value = '  +1.10%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("E25").Value = '  +2.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '666.16'
$ws.Range("E26").Value = '  +8.44%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.16'
$ws.Range("E27").Value = '  +15.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.90'
$ws.Range("E28").Value = '  +7.10%  '

$ws.Range("E29").Value = '  +4.65%  '

$ws.Range("D30").Value = '2.578.45'
$ws.Range("E30").Value = '  +2.28%  '

$ws.Range("E31").Value = '  +2.00%  '

$ws.Range("E32").Value = '  +3.56%  '

$ws.Range("E33").Value = '  +5.46%  '

$ws.Range("E34").Value = '  +4.15%  '

$ws.Range("E35").Value = '  +4.46%  '

$ws.Range("E36").Value = '  +0.22%  '

$ws.Range("E37").Value = '  +3.43%  '

$ws.Range("E38").Value = '  +3.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '153.12'
$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("E40").Value = '  +0.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.87'
$ws.Range("E41").Value = '  +2.58%  '

$ws.Range("D42").Value = '0.0₆0351'
$ws.Range("E42").Value = '  +23.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.73'
$ws.Range("E43").Value = '  +6.34%  '

$ws.Range("E44").Value = '  +3.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.32'
$ws.Range("E45").Value = '  +1.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.13'
$ws.Range("E47").Value = '  +27.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '146.63'
$ws.Range("E48").Value = '  +2.72%  '

$ws.Range("E49").Value = '  +2.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.73'
$ws.Range("E50").Value = '  +4.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.607'
$ws.Range("E51").Value = '  +2.03%  '
